$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductCreated-Event")

# Rename the MessageType value on row 3 (K3) from the old
# "ProtoBuffMessageType" label to the corrected "ProtobufType" label.
$ws.Range("K3").Value = "ProtobufType"
